$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1279.75
$ws.Range("J17").Value = 1350
$ws.Range("L17").Value = 4050
$ws.Range("N17").Value = -4386
$ws.Range("H32").Value = 1671
$ws.Range("I32").Value = 1707.8
$ws.Range("K32").Value = 1707.8
$ws.Range("M32").Value = -1381.8
$ws.Range("H40").Value = 7952.3335
$ws.Range("J40").Value = 8363.909
$ws.Range("L40").Value = 8363.909
$ws.Range("N40").Value = -8713.909
$ws.Range("H43").Value = 7875
$ws.Range("J43").Value = 8850.200000000001
$ws.Range("L43").Value = 8850.200000000001
$ws.Range("N43").Value = -8988.200000000001
$ws.Range("H53").Value = 506.9375
$ws.Range("I53").Value = 231
$ws.Range("K53").Value = 231
$ws.Range("M53").Value = 406
$ws.Range("H61").Value = 1162.5
$ws.Range("I61").Value = 1162.5
$ws.Range("K61").Value = 3487.5
$ws.Range("M61").Value = -3315.5
$ws.Range("H64").Value = 11499.286
$ws.Range("I64").Value = 8499.6
$ws.Range("J64").Value = 18998.5
$ws.Range("K64").Value = 8499.6
$ws.Range("L64").Value = 18998.5
$ws.Range("M64").Value = -8251.6
$ws.Range("N64").Value = -19494.5
$ws.Range("H67").Value = 11499.286
$ws.Range("I67").Value = 8499.6
$ws.Range("J67").Value = 18998.5
$ws.Range("K67").Value = 8499.6
$ws.Range("L67").Value = 18998.5
$ws.Range("M67").Value = -7641.6
$ws.Range("N67").Value = -20714.5
$ws.Range("H70").Value = 2585.1428
$ws.Range("I70").Value = 2182.6667
$ws.Range("K70").Value = 6548.000100000001
$ws.Range("M70").Value = -6278.000100000001
$ws.Range("H73").Value = 2585.1428
$ws.Range("I73").Value = 2182.6667
$ws.Range("K73").Value = 6548.000100000001
$ws.Range("M73").Value = -5612.000100000001
$ws.Range("H92").Value = 373.72726
$ws.Range("I92").Value = 392.1
$ws.Range("J92").Value = 190
$ws.Range("K92").Value = 392.1
$ws.Range("L92").Value = 190
$ws.Range("M92").Value = 855.9
$ws.Range("N92").Value = -2686
$ws.Range("H100").Value = 2750.6191
$ws.Range("I100").Value = 2777.0527
$ws.Range("K100").Value = 2777.0527
$ws.Range("M100").Value = -2236.0527
$ws.Range("H113").Value = 7346.6
$ws.Range("I113").Value = 2429.6667
$ws.Range("K113").Value = 2429.6667
$ws.Range("M113").Value = 824.3332999999998
$ws.Range("H131").Value = 402.5
$ws.Range("I131").Value = 402.5
$ws.Range("K131").Value = 1207.5
$ws.Range("M131").Value = 3832.5
$ws.Range("H132").Value = 2670.5
$ws.Range("I132").Value = 2448.9167
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7346.750100000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4816.750100000001
$ws.Range("N132").Value = -17060
$ws.Range("H138").Value = 10602.452
$ws.Range("I138").Value = 8485.375
$ws.Range("J138").Value = 11100.588
$ws.Range("K138").Value = 25456.125
$ws.Range("L138").Value = 33301.764
$ws.Range("M138").Value = -20316.125
$ws.Range("N138").Value = -43581.764

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7649.079
$ws.Range("I32").Value = 6685.222
$ws.Range("K32").Value = 6685.222
$ws.Range("M32").Value = -6398.222
$ws.Range("H39").Value = 396.5
$ws.Range("I39").Value = 396.5
$ws.Range("K39").Value = 396.5
$ws.Range("M39").Value = 123.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 600
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 600
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H99").Value = 2415.9092
$ws.Range("I99").Value = 2603.125
$ws.Range("K99").Value = 2603.125
$ws.Range("M99").Value = -1105.125
$ws.Range("H134").Value = 2637
$ws.Range("I134").Value = 2637
$ws.Range("K134").Value = 7911
$ws.Range("M134").Value = -5376

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4950
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4950
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H141").Value = 536812.1
$ws.Range("J141").Value = 536812.1
$ws.Range("L141").Value = 536812.1
$ws.Range("N141").Value = -547172.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3542.1428
$ws.Range("I3").Value = 3542.1428
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10626.4284
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -10514.4284
$ws.Range("N3").ClearContents()
$ws.Range("H12").Value = 877.9524
$ws.Range("J12").Value = 1216.2142
$ws.Range("L12").Value = 3648.6426
$ws.Range("N12").Value = -3994.6426
$ws.Range("H55").Value = 3399.6
$ws.Range("J55").Value = 3399.6
$ws.Range("L55").Value = 10198.8
$ws.Range("N55").Value = -10552.8
$ws.Range("H134").Value = 83336610
$ws.Range("I134").Value = 83336610
$ws.Range("K134").Value = 250009830
$ws.Range("M134").Value = -250004760

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 40000
$ws.Range("J96").Value = 40000
$ws.Range("L96").Value = 40000
$ws.Range("N96").Value = -45492
$ws.Range("H122").Value = 3601.75
$ws.Range("I122").Value = 3724.4443
$ws.Range("K122").Value = 11173.3329
$ws.Range("M122").Value = -8723.332900000001
$ws.Range("H132").Value = 3640.9
$ws.Range("I132").Value = 2301.25
$ws.Range("K132").Value = 6903.75
$ws.Range("M132").Value = -4373.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H55").Value = 994.9375
$ws.Range("I55").Value = 385.6
$ws.Range("J55").Value = 1271.909
$ws.Range("K55").Value = 385.6
$ws.Range("L55").Value = 1271.909
$ws.Range("M55").Value = -212.6
$ws.Range("N55").Value = -1617.909
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 9712.833000000001
$ws.Range("J132").Value = 9694.25
$ws.Range("L132").Value = 29082.75
$ws.Range("N132").Value = -34142.75
$ws.Range("H136").Value = 4142.3335
$ws.Range("I136").Value = 3870.8
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 11612.4
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -9062.400000000001
$ws.Range("N136").Value = -21600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H113").Value = 1205
$ws.Range("I113").Value = 730
$ws.Range("K113").Value = 2190
$ws.Range("M113").Value = -20
